# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Cilantro"
# at row 301, pushing the existing rows 301-311 down to 302-312.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 301 (existing row 301 and everything below shifts down by one,
# which is exactly the A302:R312 <- A301:R311 shift seen in the diff).
$ws.Rows.Item(301).Insert()

# Populate the newly inserted row 301 with the new weekly data point.
$ws.Cells.Item(301, 1).Value  = 4
$ws.Cells.Item(301, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(301, 3).Value  = "Los Lagos"
$ws.Cells.Item(301, 4).Value  = 44747
$ws.Cells.Item(301, 5).Value  = 10
$ws.Cells.Item(301, 6).Value  = 100112040
$ws.Cells.Item(301, 7).Value  = "Cilantro"
$ws.Cells.Item(301, 8).Value  = "Sin especificar"
$ws.Cells.Item(301, 9).Value  = "Primera"
$ws.Cells.Item(301, 10).Value = 180
$ws.Cells.Item(301, 11).Value = 12000
$ws.Cells.Item(301, 12).Value = 13000
$ws.Cells.Item(301, 13).Value = 12500
$ws.Cells.Item(301, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(301, 15).Value = "Región Metropolitana"
$ws.Cells.Item(301, 16).Value = 347
$ws.Cells.Item(301, 17).Value = 36
$ws.Cells.Item(301, 18).Value = "Hortaliza"
